$d = $word.ActiveDocument

# 1. Merge the "Specification " / "By" / " Example" runs (with proofErr
#    wrappers around "By") into a single run reading "Specification By Example".
$d.Content.Find.Execute("Specification By Example", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Specification By Example", 2) | Out-Null

# 2. Append five new bold paragraphs (list of examples) at the very end of
#    the document, right after the trailing empty paragraph and before the
#    section properties.
$texts = @(
    "Customer Not Registered For e-billing",
    "Customer Required to Update e-billing Account",
    "Billing Company Site Down",
    "Error Page Encountered",
    "?Broken Script?"
)

$insertRange = $d.Paragraphs.Item($d.Paragraphs.Count).Range

foreach ($t in $texts) {
    $insertRange.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.Text = $t
    $newPara.Range.Font.Bold = $true
    $insertRange = $newPara.Range
}
